# Support extended #ignore functionality in progbooks.
#
# The "Flows" sheet already had a "#ignore" / "Comment" block at A13:K13
# used to tell the parser to skip a row. This change extends that same
# "#ignore" / "Comment" annotation so it can also be attached per-quantity
# (next to the Birth rate and Infectiousness blocks, rows 4-5), and tidies
# up the now-unused column K underneath the Infectiousness / #ignore
# blocks. The "#ignore" cells get Excel's built-in "Neutral" cell style
# and the "Comment" cells get the built-in "Good" cell style, matching the
# existing look used elsewhere in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Flows")

# --- Birth rate block (rows 4-5): add a "#ignore" / "Comment" annotation ---
$ws.Range("H4").Value = "#ignore"
$ws.Range("H4").Style = "Neutral"

$ws.Range("I4:J4").Style = "Good"

$ws.Range("H5").Value = "#ignore"
$ws.Range("H5").Style = "Neutral"

$ws.Range("I5:J5").Style = "Good"

# --- Infectiousness block (row 11): restyle existing "Comment" cells and ---
# --- drop the now-unused column K                                        ---
$ws.Range("H11:J11").Style = "Good"
$ws.Range("K11").Clear()

# --- Existing #ignore / Comment block (row 13): restyle + drop column K ---
$ws.Range("A13").Style = "Neutral"
$ws.Range("B13:J13").Style = "Good"
$ws.Range("K13").Clear()

# Restore the saved selection/active cell as last left by the author
$ws.Range("Q5:Q6").Select()
